$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting of the Freitag block (N4:P21) into the new Samstag block (Q4:S21)
$ws.Range("N4:P21").Copy()
$ws.Range("Q4").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Update nDays counter in row 2 (J2): 5 -> 6 (Saturday added)
$ws.Range("J2").Value = 6

# Samstag header (row4) and court number (row5)
$ws.Range("Q4").Value = "Samstag"
$ws.Range("Q5").Value = 6

# Merge the new header cells like the other day blocks
$ws.Range("Q4:S4").Merge()
$ws.Range("Q5:S5").Merge()

# Hourly grid for Samstag: opening marker "x" for 06:00-12:00 (rows 6-12)
for ($r = 6; $r -le 12; $r++) {
    $ws.Cells.Item($r, 17).Value = "x"
}

Write-Output "done"
